$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 57, shifting existing rows 57:100 down to 58:101
# (weekly data refresh - a new week's price record is prepended to the
# "Vega Modelo de Temuco - Bruselas (repollito)" block).
$ws.Rows("57").Insert()

# Populate the newly inserted row 57 with this week's record.
$ws.Range("A57").Value = 10
$ws.Range("B57").Value = "Vega Modelo de Temuco"
$ws.Range("C57").Value = "La Araucanía"
$ws.Range("D57").Value = 44762
$ws.Range("E57").Value = 9
$ws.Range("F57").Value = 100112035
$ws.Range("G57").Value = "Bruselas (repollito)"
$ws.Range("H57").Value = "Sin especificar"
$ws.Range("I57").Value = "Primera"
$ws.Range("J57").Value = 25
$ws.Range("K57").Value = 26000
$ws.Range("L57").Value = 26000
$ws.Range("M57").Value = 26000
$ws.Range("N57").Value = "$/malla 10 kilos"
$ws.Range("O57").Value = "Provincia de Quillota"
$ws.Range("P57").Value = 2600
$ws.Range("Q57").Value = 10
$ws.Range("R57").Value = "Hortaliza"
